$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Julio de 2020 a las 23:45"

# Update country statistics rows (B=Casos totales, C=Nuevos casos, D=Casos activos,
# E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes)

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4624234
$ws.Range("C4").Value = 56197
$ws.Range("D4").Value = 2267120
$ws.Range("E4").Value = 2202176
$ws.Range("G4").Value = 1098
$ws.Range("H4").Value = 154938

# Row 5 - Brasil
$ws.Range("B5").Value = 2610102
$ws.Range("C5").Value = 54584
$ws.Range("E5").Value = 731420
$ws.Range("G5").Value = 1075
$ws.Range("H5").Value = 91263

# Row 22 - Francia
$ws.Range("D22").Value = 81500
$ws.Range("E22").Value = 74819

# Row 25 - Canada
$ws.Range("B25").Value = 115657
$ws.Range("C25").Value = 187
$ws.Range("D25").Value = 100686
$ws.Range("E25").Value = 6047
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = 8924

# Row 28 - Egipto
$ws.Range("B28").Value = 93757
$ws.Range("C28").Value = 401
$ws.Range("D28").Value = 38236
$ws.Range("E28").Value = 50747
$ws.Range("G28").Value = 46
$ws.Range("H28").Value = 4774

# Row 36 - Israel
$ws.Range("D36").Value = 43489
$ws.Range("E36").Value = 26047

# Row 52 - Barein
$ws.Range("B52").Value = 40755
$ws.Range("C52").Value = 444
$ws.Range("D52").Value = 37357
$ws.Range("E52").Value = 3252

# Row 63 - Moldavia
$ws.Range("B63").Value = 24343
$ws.Range("C63").Value = 396
$ws.Range("E63").Value = 6532

# Row 76 - Costa de Marfil
$ws.Range("B76").Value = 15978
$ws.Range("C76").Value = 165
$ws.Range("D76").Value = 11160
$ws.Range("E76").Value = 4718
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 100

# Row 82 - Bulgaria
$ws.Range("B82").Value = 11420
$ws.Range("C82").Value = 265
$ws.Range("D82").Value = 6173
$ws.Range("E82").Value = 4873
$ws.Range("G82").Value = 6
$ws.Range("H82").Value = 374

# Row 120 - Suazilandia
$ws.Range("B120").Value = 2577
$ws.Range("C120").Value = 26
$ws.Range("D120").Value = 1134
$ws.Range("E120").Value = 1403

# Row 129 - Ruanda
$ws.Range("B129").Value = 1994
$ws.Range("C129").Value = 31
$ws.Range("D129").Value = 1085
$ws.Range("E129").Value = 904

# Row 187 - Barbados
$ws.Range("D187").Value = 96
$ws.Range("E187").Value = 7
